$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 12749.75
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 12749.75
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 12749.75
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -13335.75
$ws.Range("H17").Value = 966.2778
$ws.Range("J17").Value = 966.2778
$ws.Range("L17").Value = 2898.8334
$ws.Range("N17").Value = -3234.8334
$ws.Range("H32").Value = 16719166
$ws.Range("I32").Value = 41666.332
$ws.Range("K32").Value = 41666.332
$ws.Range("M32").Value = -41340.332
$ws.Range("H70").Value = 2124.7144
$ws.Range("I70").Value = 3079.8
$ws.Range("K70").Value = 9239.400000000001
$ws.Range("M70").Value = -8969.400000000001
$ws.Range("H73").Value = 2124.7144
$ws.Range("I73").Value = 3079.8
$ws.Range("K73").Value = 9239.400000000001
$ws.Range("M73").Value = -8303.400000000001
$ws.Range("H76").Value = 3186.8823
$ws.Range("I76").Value = 2848.4
$ws.Range("J76").Value = 3670.4285
$ws.Range("K76").Value = 2848.4
$ws.Range("L76").Value = 3670.4285
$ws.Range("M76").Value = -2533.4
$ws.Range("N76").Value = -4300.4285
$ws.Range("H79").Value = 3186.8823
$ws.Range("I79").Value = 2848.4
$ws.Range("J79").Value = 3670.4285
$ws.Range("K79").Value = 2848.4
$ws.Range("L79").Value = 3670.4285
$ws.Range("M79").Value = -1756.4
$ws.Range("N79").Value = -5854.4285
$ws.Range("H86").Value = 3335.125
$ws.Range("I86").Value = 2804.5715
$ws.Range("J86").Value = 4077.9
$ws.Range("K86").Value = 2804.5715
$ws.Range("L86").Value = 4077.9
$ws.Range("M86").Value = -1681.5715
$ws.Range("N86").Value = -6323.9
$ws.Range("H89").Value = 3335.125
$ws.Range("I89").Value = 2804.5715
$ws.Range("J89").Value = 4077.9
$ws.Range("K89").Value = 14022.8575
$ws.Range("L89").Value = 20389.5
$ws.Range("M89").Value = -8406.8575
$ws.Range("N89").Value = -31621.5
$ws.Range("H97").Value = 937.8333
$ws.Range("J97").Value = 937.8333
$ws.Range("L97").Value = 2813.4999
$ws.Range("N97").Value = -3805.4999
$ws.Range("H98").Value = 638.8
$ws.Range("I98").Value = 644.625
$ws.Range("K98").Value = 644.625
$ws.Range("M98").Value = 853.375
$ws.Range("H100").Value = 1572
$ws.Range("I100").Value = 1286.3334
$ws.Range("K100").Value = 1286.3334
$ws.Range("M100").Value = -745.3334
$ws.Range("H106").Value = 4179.8887
$ws.Range("I106").Value = 4018.6924
$ws.Range("J106").Value = 4599
$ws.Range("K106").Value = 4018.6924
$ws.Range("L106").Value = 4599
$ws.Range("M106").Value = -3387.6924
$ws.Range("N106").Value = -5861
$ws.Range("H107").Value = 447.75
$ws.Range("I107").Value = 498
$ws.Range("J107").Value = 397.5
$ws.Range("K107").Value = 498
$ws.Range("L107").Value = 397.5
$ws.Range("M107").Value = 1422
$ws.Range("N107").Value = -4237.5
$ws.Range("H111").Value = 702.8182
$ws.Range("I111").Value = 518.2
$ws.Range("J111").Value = 856.6667
$ws.Range("K111").Value = 1554.6
$ws.Range("L111").Value = 2570.0001
$ws.Range("M111").Value = 1512.4
$ws.Range("N111").Value = -8704.000100000001
$ws.Range("H112").Value = 887.53845
$ws.Range("I112").Value = 634.5
$ws.Range("J112").Value = 1000
$ws.Range("K112").Value = 1903.5
$ws.Range("L112").Value = 3000
$ws.Range("M112").Value = -795.5
$ws.Range("N112").Value = -5216
$ws.Range("H122").Value = 638.8
$ws.Range("I122").Value = 644.625
$ws.Range("K122").Value = 1933.875
$ws.Range("M122").Value = 516.125
$ws.Range("H130").Value = 60780
$ws.Range("J130").Value = 60780
$ws.Range("L130").Value = 60780
$ws.Range("N130").Value = -70820
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1575.375
$ws.Range("I132").Value = 1575.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4726.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2196.125
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 469869.9
$ws.Range("I137").Value = 1407.8636
$ws.Range("J137").Value = 1614999.4
$ws.Range("K137").Value = 4223.5908
$ws.Range("L137").Value = 4844998.199999999
$ws.Range("M137").Value = -1673.5908
$ws.Range("N137").Value = -4850098.199999999
$ws.Range("H138").Value = 1626.661
$ws.Range("I138").Value = 1190.0322
$ws.Range("J138").Value = 2110.0715
$ws.Range("K138").Value = 3570.0966
$ws.Range("L138").Value = 6330.2145
$ws.Range("M138").Value = 1569.9034
$ws.Range("N138").Value = -16610.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 406.8889
$ws.Range("I4").Value = 168.25
$ws.Range("K4").Value = 168.25
$ws.Range("M4").Value = -52.25
$ws.Range("H45").Value = 17860844
$ws.Range("I45").Value = 4206
$ws.Range("J45").Value = 25003500
$ws.Range("K45").Value = 4206
$ws.Range("L45").Value = 25003500
$ws.Range("M45").Value = -3829
$ws.Range("N45").Value = -25004254
$ws.Range("H61").Value = 47862.637
$ws.Range("I61").Value = 2234.7646
$ws.Range("J61").Value = 202997.4
$ws.Range("K61").Value = 2234.7646
$ws.Range("L61").Value = 202997.4
$ws.Range("M61").Value = -2022.7646
$ws.Range("N61").Value = -203421.4
$ws.Range("H74").Value = 33122.844
$ws.Range("I74").Value = 38253
$ws.Range("J74").Value = 5420
$ws.Range("K74").Value = 38253
$ws.Range("L74").Value = 5420
$ws.Range("M74").Value = -37379
$ws.Range("N74").Value = -7168
$ws.Range("H77").Value = 33122.844
$ws.Range("I77").Value = 38253
$ws.Range("J77").Value = 5420
$ws.Range("K77").Value = 191265
$ws.Range("L77").Value = 27100
$ws.Range("M77").Value = -186897
$ws.Range("N77").Value = -35836
$ws.Range("H97").Value = 958.4
$ws.Range("I97").Value = 1059.7693
$ws.Range("J97").Value = 299.5
$ws.Range("K97").Value = 1059.7693
$ws.Range("L97").Value = 299.5
$ws.Range("M97").Value = -563.7692999999999
$ws.Range("N97").Value = -1291.5
$ws.Range("H102").Value = 76082.664
$ws.Range("I102").Value = 112583.78
$ws.Range("J102").Value = 21331
$ws.Range("K102").Value = 112583.78
$ws.Range("L102").Value = 21331
$ws.Range("M102").Value = -110961.78
$ws.Range("N102").Value = -24575
$ws.Range("H110").Value = 2549.8
$ws.Range("I110").Value = 2249.5
$ws.Range("J110").Value = 2750
$ws.Range("K110").Value = 2249.5
$ws.Range("L110").Value = 2750
$ws.Range("M110").Value = -204.5
$ws.Range("N110").Value = -6840
$ws.Range("H122").Value = 4999.25
$ws.Range("I122").Value = 4999
$ws.Range("K122").Value = 14997
$ws.Range("M122").Value = -12547
$ws.Range("H132").Value = 2436.2856
$ws.Range("I132").Value = 1983.4783
$ws.Range("J132").Value = 4519.2
$ws.Range("K132").Value = 5950.4349
$ws.Range("L132").Value = 13557.6
$ws.Range("M132").Value = -3420.4349
$ws.Range("N132").Value = -18617.6
$ws.Range("H136").Value = 47862.637
$ws.Range("I136").Value = 2234.7646
$ws.Range("J136").Value = 202997.4
$ws.Range("K136").Value = 6704.293799999999
$ws.Range("L136").Value = 608992.2
$ws.Range("M136").Value = -4154.293799999999
$ws.Range("N136").Value = -614092.2
$ws.Range("H139").Value = 124857.5
$ws.Range("J139").Value = 124857.5
$ws.Range("L139").Value = 124857.5
$ws.Range("N139").Value = -135137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H64").Value = 1572.3334
$ws.Range("I64").Value = 1171.6666
$ws.Range("J64").Value = 1772.6666
$ws.Range("K64").Value = 1171.6666
$ws.Range("L64").Value = 1772.6666
$ws.Range("M64").Value = -946.6666
$ws.Range("N64").Value = -2222.6666
$ws.Range("H67").Value = 1572.3334
$ws.Range("I67").Value = 1171.6666
$ws.Range("J67").Value = 1772.6666
$ws.Range("K67").Value = 1171.6666
$ws.Range("L67").Value = 1772.6666
$ws.Range("M67").Value = -391.6666
$ws.Range("N67").Value = -3332.6666
$ws.Range("H86").Value = 3316.4814
$ws.Range("I86").Value = 3570.6667
$ws.Range("J86").Value = 2998.75
$ws.Range("K86").Value = 3570.6667
$ws.Range("L86").Value = 2998.75
$ws.Range("M86").Value = -2447.6667
$ws.Range("N86").Value = -5244.75
$ws.Range("H89").Value = 3316.4814
$ws.Range("I89").Value = 3570.6667
$ws.Range("J89").Value = 2998.75
$ws.Range("K89").Value = 17853.3335
$ws.Range("L89").Value = 14993.75
$ws.Range("M89").Value = -12237.3335
$ws.Range("N89").Value = -26225.75
$ws.Range("H105").Value = 129374.625
$ws.Range("I105").Value = 336336.66
$ws.Range("J105").Value = 5197.4
$ws.Range("K105").Value = 336336.66
$ws.Range("L105").Value = 5197.4
$ws.Range("M105").Value = -334589.66
$ws.Range("N105").Value = -8691.4
$ws.Range("H107").Value = 5885175.5
$ws.Range("I107").Value = 7695091
$ws.Range("J107").Value = 2949.75
$ws.Range("K107").Value = 7695091
$ws.Range("L107").Value = 2949.75
$ws.Range("M107").Value = -7693171
$ws.Range("N107").Value = -6789.75
$ws.Range("H132").Value = 28160.285
$ws.Range("J132").Value = 28160.285
$ws.Range("L132").Value = 28160.285
$ws.Range("N132").Value = -38280.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6474.9375
$ws.Range("I7").Value = 6050
$ws.Range("K7").Value = 6050
$ws.Range("M7").Value = -5937
$ws.Range("H16").Value = 1750.1428
$ws.Range("I16").Value = 1291.7142
$ws.Range("J16").Value = 2208.5715
$ws.Range("K16").Value = 1291.7142
$ws.Range("L16").Value = 2208.5715
$ws.Range("M16").Value = -1004.7142
$ws.Range("N16").Value = -2782.5715
$ws.Range("H22").Value = 320.91666
$ws.Range("I22").Value = 328.27274
$ws.Range("J22").Value = 240
$ws.Range("K22").Value = 328.27274
$ws.Range("L22").Value = 240
$ws.Range("M22").Value = 21.72726
$ws.Range("N22").Value = -940
$ws.Range("H31").Value = 3613.5
$ws.Range("I31").Value = 2085.5
$ws.Range("J31").Value = 5905.5
$ws.Range("K31").Value = 2085.5
$ws.Range("L31").Value = 5905.5
$ws.Range("M31").Value = -1790.5
$ws.Range("N31").Value = -6495.5
$ws.Range("H34").Value = 3613.5
$ws.Range("I34").Value = 2085.5
$ws.Range("J34").Value = 5905.5
$ws.Range("K34").Value = 2085.5
$ws.Range("L34").Value = 5905.5
$ws.Range("M34").Value = -1883.5
$ws.Range("N34").Value = -6309.5
$ws.Range("H58").Value = 1362.5
$ws.Range("I58").Value = 450
$ws.Range("J58").Value = 2275
$ws.Range("K58").Value = 450
$ws.Range("L58").Value = 2275
$ws.Range("M58").Value = -247
$ws.Range("N58").Value = -2681
$ws.Range("H74").Value = 56333
$ws.Range("J74").Value = 56333
$ws.Range("L74").Value = 56333
$ws.Range("N74").Value = -58081
$ws.Range("H77").Value = 56333
$ws.Range("J77").Value = 56333
$ws.Range("L77").Value = 168999
$ws.Range("N77").Value = -177735
$ws.Range("H107").Value = 1400.3572
$ws.Range("J107").Value = 1438.2858
$ws.Range("L107").Value = 1438.2858
$ws.Range("N107").Value = -5278.2858
$ws.Range("H113").Value = 1750.1428
$ws.Range("I113").Value = 1291.7142
$ws.Range("J113").Value = 2208.5715
$ws.Range("K113").Value = 1291.7142
$ws.Range("L113").Value = 2208.5715
$ws.Range("M113").Value = 878.2858000000001
$ws.Range("N113").Value = -6548.5715
$ws.Range("H132").Value = 2924957.5
$ws.Range("I132").Value = 2843102
$ws.Range("J132").Value = 3252378.5
$ws.Range("K132").Value = 8529306
$ws.Range("L132").Value = 9757135.5
$ws.Range("M132").Value = -8526776
$ws.Range("N132").Value = -9762195.5
$ws.Range("H134").Value = 3815097.5
$ws.Range("I134").Value = 5105382.5
$ws.Range("K134").Value = 15316147.5
$ws.Range("M134").Value = -15313612.5
$ws.Range("H136").Value = 1362.5
$ws.Range("I136").Value = 450
$ws.Range("J136").Value = 2275
$ws.Range("K136").Value = 1350
$ws.Range("L136").Value = 6825
$ws.Range("M136").Value = 1200
$ws.Range("N136").Value = -11925

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 750951.6
$ws.Range("I22").Value = 791
$ws.Range("J22").Value = 1501112.2
$ws.Range("K22").Value = 2373
$ws.Range("L22").Value = 4503336.6
$ws.Range("M22").Value = -2204
$ws.Range("N22").Value = -4503674.6
$ws.Range("H23").Value = 79021.234
$ws.Range("J23").Value = 85442.336
$ws.Range("L23").Value = 256327.008
$ws.Range("N23").Value = -256797.008
$ws.Range("H25").Value = 229
$ws.Range("I25").Value = 229
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 687
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -518
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 750951.6
$ws.Range("I27").Value = 791
$ws.Range("J27").Value = 1501112.2
$ws.Range("K27").Value = 2373
$ws.Range("L27").Value = 4503336.6
$ws.Range("M27").Value = -2271
$ws.Range("N27").Value = -4503540.6
$ws.Range("H29").Value = 33333646
$ws.Range("I29").Value = 101.333336
$ws.Range("J29").Value = 50000416
$ws.Range("K29").Value = 304.000008
$ws.Range("L29").Value = 150001248
$ws.Range("M29").Value = -27.00000799999998
$ws.Range("N29").Value = -150001802
$ws.Range("H30").Value = 229
$ws.Range("I30").Value = 229
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 687
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -585
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 1973.75
$ws.Range("I32").Value = 263.33334
$ws.Range("K32").Value = 790.0000200000001
$ws.Range("M32").Value = -507.0000200000001
$ws.Range("H80").Value = 1598
$ws.Range("I80").Value = 1598
$ws.Range("K80").Value = 4794
$ws.Range("M80").Value = -3858
$ws.Range("H83").Value = 1598
$ws.Range("I83").Value = 1598
$ws.Range("K83").Value = 14382
$ws.Range("M83").Value = -9702
$ws.Range("H110").Value = 7359.8
$ws.Range("I110").Value = 7399.5
$ws.Range("J110").Value = 7333.3335
$ws.Range("K110").Value = 22198.5
$ws.Range("L110").Value = 22000.0005
$ws.Range("M110").Value = -18108.5
$ws.Range("N110").Value = -30180.0005
$ws.Range("H131").Value = 1308
$ws.Range("I131").Value = 949.55554
$ws.Range("J131").Value = 2383.3333
$ws.Range("K131").Value = 2848.66662
$ws.Range("L131").Value = 7149.999899999999
$ws.Range("M131").Value = 2191.33338
$ws.Range("N131").Value = -17229.9999
$ws.Range("H137").Value = 3583
$ws.Range("J137").Value = 4959.1
$ws.Range("L137").Value = 14877.3
$ws.Range("N137").Value = -25077.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11882750
$ws.Range("I14").Value = 11882750
$ws.Range("K14").Value = 11882750
$ws.Range("M14").Value = -11882582
$ws.Range("H75").Value = 10000
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 10000
$ws.Range("N75").Value = -11748
$ws.Range("H78").Value = 10000
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 30000
$ws.Range("N78").Value = -38736
$ws.Range("H80").Value = 3739.8
$ws.Range("I80").Value = 5999
$ws.Range("J80").Value = 3175
$ws.Range("K80").Value = 5999
$ws.Range("L80").Value = 3175
$ws.Range("M80").Value = -5001
$ws.Range("N80").Value = -5171
$ws.Range("H83").Value = 3739.8
$ws.Range("I83").Value = 5999
$ws.Range("J83").Value = 3175
$ws.Range("K83").Value = 29995
$ws.Range("L83").Value = 15875
$ws.Range("M83").Value = -25003
$ws.Range("N83").Value = -25859
$ws.Range("H132").Value = 10004
$ws.Range("I132").Value = 10004
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30012
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -27482
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3353333.2
$ws.Range("J20").Value = 3353333.2
$ws.Range("L20").Value = 3353333.2
$ws.Range("N20").Value = -3353785.2
$ws.Range("H22").Value = 9623.5
$ws.Range("I22").Value = 1407.4546
$ws.Range("K22").Value = 1407.4546
$ws.Range("M22").Value = -1112.4546
$ws.Range("H27").Value = 9623.5
$ws.Range("I27").Value = 1407.4546
$ws.Range("K27").Value = 1407.4546
$ws.Range("M27").Value = -1300.4546
$ws.Range("H61").Value = 715.38464
$ws.Range("I61").Value = 307
$ws.Range("K61").Value = 307
$ws.Range("M61").Value = -105
$ws.Range("H63").Value = 27563.75
$ws.Range("J63").Value = 30085
$ws.Range("L63").Value = 30085
$ws.Range("N63").Value = -31583
$ws.Range("H64").Value = 16000
$ws.Range("J64").Value = 16000
$ws.Range("L64").Value = 16000
$ws.Range("N64").Value = -16450
$ws.Range("H66").Value = 27563.75
$ws.Range("J66").Value = 30085
$ws.Range("L66").Value = 90255
$ws.Range("N66").Value = -97743
$ws.Range("H67").Value = 16000
$ws.Range("J67").Value = 16000
$ws.Range("L67").Value = 16000
$ws.Range("N67").Value = -17560
$ws.Range("H68").Value = 165519
$ws.Range("I68").Value = 237310.78
$ws.Range("J68").Value = 3987.5
$ws.Range("K68").Value = 237310.78
$ws.Range("L68").Value = 3987.5
$ws.Range("M68").Value = -236561.78
$ws.Range("N68").Value = -5485.5
$ws.Range("H71").Value = 165519
$ws.Range("I71").Value = 237310.78
$ws.Range("J71").Value = 3987.5
$ws.Range("K71").Value = 1186553.9
$ws.Range("L71").Value = 19937.5
$ws.Range("M71").Value = -1182809.9
$ws.Range("N71").Value = -27425.5
$ws.Range("H82").Value = 2400.25
$ws.Range("I82").Value = 2200.3333
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 2200.3333
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1839.3333
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2400.25
$ws.Range("I85").Value = 2200.3333
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 2200.3333
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -952.3332999999998
$ws.Range("N85").Value = -5496
$ws.Range("H93").Value = 2108.8
$ws.Range("I93").Value = 1898.8889
$ws.Range("K93").Value = 1898.8889
$ws.Range("M93").Value = -650.8888999999999
$ws.Range("H113").Value = 715.38464
$ws.Range("I113").Value = 307
$ws.Range("K113").Value = 307
$ws.Range("M113").Value = 1863
$ws.Range("H132").Value = 3339.9443
$ws.Range("I132").Value = 2540.625
$ws.Range("K132").Value = 7621.875
$ws.Range("M132").Value = -5091.875
$ws.Range("H136").Value = 1871.6
$ws.Range("I136").Value = 1571.9445
$ws.Range("J136").Value = 2642.1428
$ws.Range("K136").Value = 4715.833500000001
$ws.Range("L136").Value = 7926.428400000001
$ws.Range("M136").Value = -2165.833500000001
$ws.Range("N136").Value = -13026.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H8").Value = 9999
$ws.Range("J8").Value = 9999
$ws.Range("L8").Value = 9999
$ws.Range("N8").Value = -10279
$ws.Range("H13").Value = 1500
$ws.Range("J13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("N13").Value = -1780
$ws.Range("H54").Value = 52999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 52999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 52999
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -54039
$ws.Range("H62").Value = 5538
$ws.Range("I62").Value = 5538
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5538
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4914
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5538
$ws.Range("I65").Value = 5538
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 27690
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -24570
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1313
$ws.Range("I113").Value = 1437.4
$ws.Range("J113").Value = 898.3333
$ws.Range("K113").Value = 4312.200000000001
$ws.Range("L113").Value = 2694.9999
$ws.Range("M113").Value = -2142.200000000001
$ws.Range("N113").Value = -7034.9999
$ws.Range("H132").Value = 1990.3793
$ws.Range("I132").Value = 1835.4231
$ws.Range("K132").Value = 5506.2693
$ws.Range("M132").Value = -2976.2693
$ws.Range("H136").Value = 1433.1818
$ws.Range("I136").Value = 1276.5
$ws.Range("K136").Value = 3829.5
$ws.Range("M136").Value = -1279.5
